$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor TPM-derived values (NATMI LR-pairs output)
# Columns: G/H = ligand avg/total expr, I/J = ligand specificity (avg/total)
#          M/N = receptor avg/total expr, O/P = receptor specificity (avg/total)
#          Q/R = edge avg/total expr weight, S/T = edge avg/total specificity
# Row 2
$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 10.35074139316533
$ws.Range("R2").Value = 93.156672538488
$ws.Range("S2").Value = 0.003636201619892214
$ws.Range("T2").Value = 0.003636201619892215
# Row 3
$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 843.7804190570751
$ws.Range("R3").Value = 7594.023771513675
$ws.Range("S3").Value = 0.2964189336847496
$ws.Range("T3").Value = 0.2964189336847497
# Row 4
$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 64.91393798364565
$ws.Range("R4").Value = 584.2254418528109
$ws.Range("S4").Value = 0.02280417967021899
$ws.Range("T4").Value = 0.022804179670219
# Row 5
$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 7.34014536632
$ws.Range("R5").Value = 66.06130829688
$ws.Range("S5").Value = 0.002578583258671777
$ws.Range("T5").Value = 0.002578583258671777
# Row 6
$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("Q6").Value = 598.3601268623066
$ws.Range("S6").Value = 0.2102031129878579
$ws.Range("T6").Value = 0.2102031129878579
# Row 7
$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 46.03319926579
$ws.Range("R7").Value = 414.2987933921099
$ws.Range("S7").Value = 0.01617140138865925
$ws.Range("T7").Value = 0.01617140138865926
# Row 8
$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 14.36871602438133
$ws.Range("R8").Value = 129.318444219432
$ws.Range("S8").Value = 0.005047710738684481
$ws.Range("T8").Value = 0.005047710738684482
# Row 9
$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 1171.321045308907
$ws.Range("R9").Value = 10541.88940778016
$ws.Range("S9").Value = 0.4114835179998258
$ws.Range("T9").Value = 0.4114835179998259
# Row 10
$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 90.11237992354765
$ws.Range("R10").Value = 811.0114193119289
$ws.Range("S10").Value = 0.03165635865143992
$ws.Range("T10").Value = 0.03165635865143992
